# Adding Master Data XLS
# The "Vostro" / DKS biometric-device rows (id 589 and 638) are stale test
# data and get removed entirely; every row below them shifts up by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows for the Vostro/DKS device (old rows 6 & 7) - this
# shifts rows 8-13 up to become the new rows 6-11.
$ws.Rows("6:7").Delete()

# Leave the selection where the author's last interaction landed.
$ws.Range("E16").Select() | Out-Null

# Touch the page setup (paper size / orientation) as in the final workbook.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
